# Update the "Estado de Cuenta" worker table: the record for
# LEONARDO JAVIER VERA DORIA (doc 1047482742, periodo 1708) now comes
# before the record for LORENA BEATRIZ DONADO LOPEZ (doc 45757837,
# periodo 1712) - i.e. rows 17 and 18 swap their identifying data while
# "Valor Mora"/"Salario Basico" (columns F:G) stay as they were.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = "1047482742"
$ws.Range("D17").Value = "LEONARDO JAVIER VERA DORIA"
$ws.Range("E17").Value = "1708"

$ws.Range("C18").Value = "45757837"
$ws.Range("D18").Value = "LORENA BEATRIZ DONADO LOPEZ"
$ws.Range("E18").Value = "1712"
